# Auto-generated script applying scheduled market-price update to Maduin_Profits workbook
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(3, 8).Value = 30000
$ws.Cells.Item(3, 10).Value = 30000
$ws.Cells.Item(3, 12).Value = 30000
$ws.Cells.Item(3, 14).Value = -30228

$ws.Cells.Item(6, 8).Value = 202.14285
$ws.Cells.Item(6, 9).Value = 219.16667
$ws.Cells.Item(6, 11).Value = 657.50001
$ws.Cells.Item(6, 13).Value = -545.50001

$ws.Cells.Item(8, 8).Value = 64.44444
$ws.Cells.Item(8, 9).Value = 60
$ws.Cells.Item(8, 11).Value = 180
$ws.Cells.Item(8, 13).Value = -41

$ws.Cells.Item(31, 8).Value = 8
$ws.Cells.Item(31, 9).Value = 8
$ws.Cells.Item(31, 10).Value = 0
$ws.Cells.Item(31, 11).Value = 24
$ws.Cells.Item(31, 12).Value = 0
$ws.Cells.Item(31, 13).Value = 206
$ws.Cells.Item(31, 14).ClearContents()

$ws.Cells.Item(39, 8).Value = 454.1111
$ws.Cells.Item(39, 9).Value = 121
$ws.Cells.Item(39, 10).Value = 720.6
$ws.Cells.Item(39, 11).Value = 363
$ws.Cells.Item(39, 12).Value = 2161.8
$ws.Cells.Item(39, 13).Value = -67
$ws.Cells.Item(39, 14).Value = -2753.8

$ws.Cells.Item(55, 8).Value = 149.54546
$ws.Cells.Item(55, 9).Value = 99.55556
$ws.Cells.Item(55, 10).Value = 374.5
$ws.Cells.Item(55, 11).Value = 99.55556
$ws.Cells.Item(55, 12).Value = 374.5
$ws.Cells.Item(55, 13).Value = 114.44444
$ws.Cells.Item(55, 14).Value = -802.5

$ws.Cells.Item(64, 8).Value = 13499.375
$ws.Cells.Item(64, 9).Value = 8888.444
$ws.Cells.Item(64, 11).Value = 8888.444
$ws.Cells.Item(64, 13).Value = -8640.444

$ws.Cells.Item(67, 8).Value = 13499.375
$ws.Cells.Item(67, 9).Value = 8888.444
$ws.Cells.Item(67, 11).Value = 8888.444
$ws.Cells.Item(67, 13).Value = -8030.444

$ws.Cells.Item(86, 8).Value = 4666.6665
$ws.Cells.Item(86, 9).Value = 3500
$ws.Cells.Item(86, 10).Value = 7000
$ws.Cells.Item(86, 11).Value = 3500
$ws.Cells.Item(86, 12).Value = 7000
$ws.Cells.Item(86, 13).Value = -2377
$ws.Cells.Item(86, 14).Value = -9246

$ws.Cells.Item(89, 8).Value = 4666.6665
$ws.Cells.Item(89, 9).Value = 3500
$ws.Cells.Item(89, 10).Value = 7000
$ws.Cells.Item(89, 11).Value = 17500
$ws.Cells.Item(89, 12).Value = 35000
$ws.Cells.Item(89, 13).Value = -11884
$ws.Cells.Item(89, 14).Value = -46232

$ws.Cells.Item(92, 8).Value = 368.27274
$ws.Cells.Item(92, 9).Value = 387.6
$ws.Cells.Item(92, 10).Value = 175
$ws.Cells.Item(92, 11).Value = 387.6
$ws.Cells.Item(92, 12).Value = 175
$ws.Cells.Item(92, 13).Value = 860.4
$ws.Cells.Item(92, 14).Value = -2671

$ws.Cells.Item(102, 8).Value = 30000
$ws.Cells.Item(102, 10).Value = 30000
$ws.Cells.Item(102, 12).Value = 30000
$ws.Cells.Item(102, 14).Value = -36490

$ws.Cells.Item(131, 8).Value = 410
$ws.Cells.Item(131, 9).Value = 410
$ws.Cells.Item(131, 11).Value = 1230
$ws.Cells.Item(131, 13).Value = 3810

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 7776.9487
$ws.Cells.Item(32, 9).Value = 6418.9443
$ws.Cells.Item(32, 11).Value = 6418.9443
$ws.Cells.Item(32, 13).Value = -6131.9443

$ws.Cells.Item(45, 8).Value = 2554.2727
$ws.Cells.Item(45, 9).Value = 2031.5
$ws.Cells.Item(45, 11).Value = 2031.5
$ws.Cells.Item(45, 13).Value = -1654.5

$ws.Cells.Item(74, 8).Value = 2899.9
$ws.Cells.Item(74, 9).Value = 2899.9
$ws.Cells.Item(74, 11).Value = 2899.9
$ws.Cells.Item(74, 13).Value = -2025.9

$ws.Cells.Item(77, 8).Value = 2899.9
$ws.Cells.Item(77, 9).Value = 2899.9
$ws.Cells.Item(77, 11).Value = 14499.5
$ws.Cells.Item(77, 13).Value = -10131.5

$ws.Cells.Item(132, 8).Value = 2800.2
$ws.Cells.Item(132, 9).Value = 1812.75
$ws.Cells.Item(132, 11).Value = 5438.25
$ws.Cells.Item(132, 13).Value = -2908.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(94, 8).Value = 1981.1333
$ws.Cells.Item(94, 9).Value = 1362.8462
$ws.Cells.Item(94, 10).Value = 6000
$ws.Cells.Item(94, 11).Value = 1362.8462
$ws.Cells.Item(94, 12).Value = 6000
$ws.Cells.Item(94, 13).Value = -911.8462
$ws.Cells.Item(94, 14).Value = -6902

$ws.Cells.Item(99, 8).Value = 1416.3334
$ws.Cells.Item(99, 9).Value = 499
$ws.Cells.Item(99, 10).Value = 1875
$ws.Cells.Item(99, 11).Value = 499
$ws.Cells.Item(99, 12).Value = 1875
$ws.Cells.Item(99, 13).Value = 999
$ws.Cells.Item(99, 14).Value = -4871

$ws.Cells.Item(100, 8).Value = 12000
$ws.Cells.Item(100, 10).Value = 12000
$ws.Cells.Item(100, 12).Value = 12000
$ws.Cells.Item(100, 14).Value = -14164

$ws.Cells.Item(107, 8).Value = 1400
$ws.Cells.Item(107, 9).Value = 1400
$ws.Cells.Item(107, 11).Value = 1400
$ws.Cells.Item(107, 13).Value = 520

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(28, 8).Value = 16821.5
$ws.Cells.Item(28, 10).Value = 16821.5
$ws.Cells.Item(28, 12).Value = 16821.5
$ws.Cells.Item(28, 14).Value = -17311.5

$ws.Cells.Item(31, 8).Value = 3195.4614
$ws.Cells.Item(31, 9).Value = 3003.9092
$ws.Cells.Item(31, 11).Value = 3003.9092
$ws.Cells.Item(31, 13).Value = -2708.9092

$ws.Cells.Item(34, 8).Value = 3195.4614
$ws.Cells.Item(34, 9).Value = 3003.9092
$ws.Cells.Item(34, 11).Value = 3003.9092
$ws.Cells.Item(34, 13).Value = -2801.9092

$ws.Cells.Item(99, 8).Value = 9505.799999999999
$ws.Cells.Item(99, 9).Value = 8767.25
$ws.Cells.Item(99, 11).Value = 8767.25
$ws.Cells.Item(99, 13).Value = -7269.25

$ws.Cells.Item(126, 8).Value = 9505.799999999999
$ws.Cells.Item(126, 9).Value = 8767.25
$ws.Cells.Item(126, 11).Value = 26301.75
$ws.Cells.Item(126, 13).Value = -23831.75

$ws.Cells.Item(132, 8).Value = 6868.654
$ws.Cells.Item(132, 9).Value = 2639.2
$ws.Cells.Item(132, 11).Value = 7917.599999999999
$ws.Cells.Item(132, 13).Value = -5387.599999999999

$ws.Cells.Item(141, 8).Value = 296045.47
$ws.Cells.Item(141, 10).Value = 296045.47
$ws.Cells.Item(141, 12).Value = 296045.47
$ws.Cells.Item(141, 14).Value = -306405.47

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(75, 8).Value = 874.25
$ws.Cells.Item(75, 9).Value = 866.1667
$ws.Cells.Item(75, 10).Value = 898.5
$ws.Cells.Item(75, 11).Value = 2598.5001
$ws.Cells.Item(75, 12).Value = 2695.5
$ws.Cells.Item(75, 13).Value = -1600.5001
$ws.Cells.Item(75, 14).Value = -4691.5

$ws.Cells.Item(78, 8).Value = 874.25
$ws.Cells.Item(78, 9).Value = 866.1667
$ws.Cells.Item(78, 10).Value = 898.5
$ws.Cells.Item(78, 11).Value = 7795.5003
$ws.Cells.Item(78, 12).Value = 8086.5
$ws.Cells.Item(78, 13).Value = -2803.5003
$ws.Cells.Item(78, 14).Value = -18070.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(12, 8).Value = 0
$ws.Cells.Item(12, 10).Value = 0
$ws.Cells.Item(12, 12).Value = 0
$ws.Cells.Item(12, 14).ClearContents()

$ws.Cells.Item(17, 8).Value = 1400
$ws.Cells.Item(17, 10).Value = 1200
$ws.Cells.Item(17, 12).Value = 1200
$ws.Cells.Item(17, 14).Value = -1536

$ws.Cells.Item(80, 8).Value = 3236.8
$ws.Cells.Item(80, 10).Value = 2799.3333
$ws.Cells.Item(80, 12).Value = 2799.3333
$ws.Cells.Item(80, 14).Value = -4795.3333

$ws.Cells.Item(83, 8).Value = 3236.8
$ws.Cells.Item(83, 10).Value = 2799.3333
$ws.Cells.Item(83, 12).Value = 13996.6665
$ws.Cells.Item(83, 14).Value = -23980.6665

$ws.Cells.Item(92, 8).Value = 13927.8
$ws.Cells.Item(92, 10).Value = 13927.8
$ws.Cells.Item(92, 12).Value = 13927.8
$ws.Cells.Item(92, 14).Value = -17671.8

$ws.Cells.Item(98, 8).Value = 9000
$ws.Cells.Item(98, 10).Value = 9000
$ws.Cells.Item(98, 12).Value = 9000
$ws.Cells.Item(98, 14).Value = -14990

$ws.Cells.Item(99, 8).Value = 8578.166999999999
$ws.Cells.Item(99, 9).Value = 8578.166999999999
$ws.Cells.Item(99, 11).Value = 8578.166999999999
$ws.Cells.Item(99, 13).Value = -6332.166999999999

$ws.Cells.Item(102, 8).Value = 0
$ws.Cells.Item(102, 9).Value = 0
$ws.Cells.Item(102, 11).Value = 0
$ws.Cells.Item(102, 13).ClearContents()

$ws.Cells.Item(103, 8).Value = 73333.336
$ws.Cells.Item(103, 10).Value = 73333.336
$ws.Cells.Item(103, 12).Value = 73333.336
$ws.Cells.Item(103, 14).Value = -75677.336

$ws.Cells.Item(107, 8).Value = 569.44446
$ws.Cells.Item(107, 9).Value = 578.125
$ws.Cells.Item(107, 10).Value = 500
$ws.Cells.Item(107, 11).Value = 578.125
$ws.Cells.Item(107, 12).Value = 500
$ws.Cells.Item(107, 13).Value = 1341.875
$ws.Cells.Item(107, 14).Value = -4340

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(2, 8).Value = 3680

$ws.Cells.Item(46, 8).Value = 1528.1538
$ws.Cells.Item(46, 9).Value = 1562.3334
$ws.Cells.Item(46, 10).Value = 1498.8572
$ws.Cells.Item(46, 11).Value = 1562.3334
$ws.Cells.Item(46, 12).Value = 1498.8572
$ws.Cells.Item(46, 13).Value = -1374.3334
$ws.Cells.Item(46, 14).Value = -1874.8572

$ws.Cells.Item(93, 8).Value = 747.1667
$ws.Cells.Item(93, 9).Value = 820.75
$ws.Cells.Item(93, 10).Value = 600
$ws.Cells.Item(93, 11).Value = 820.75
$ws.Cells.Item(93, 12).Value = 600
$ws.Cells.Item(93, 13).Value = 427.25
$ws.Cells.Item(93, 14).Value = -3096

$ws.Cells.Item(122, 8).Value = 8090.2856
$ws.Cells.Item(122, 9).Value = 8147
$ws.Cells.Item(122, 11).Value = 24441
$ws.Cells.Item(122, 13).Value = -21991

$ws.Cells.Item(132, 8).Value = 7966.6665
$ws.Cells.Item(132, 9).Value = 6783.3335
$ws.Cells.Item(132, 10).Value = 10333.333
$ws.Cells.Item(132, 11).Value = 20350.0005
$ws.Cells.Item(132, 12).Value = 30999.999
$ws.Cells.Item(132, 13).Value = -17820.0005
$ws.Cells.Item(132, 14).Value = -36059.999

$ws.Cells.Item(136, 8).Value = 3811.4443
$ws.Cells.Item(136, 9).Value = 3329
$ws.Cells.Item(136, 11).Value = 9987
$ws.Cells.Item(136, 13).Value = -7437

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(56, 8).Value = 29999.5
$ws.Cells.Item(56, 9).Value = 10000
$ws.Cells.Item(56, 10).Value = 49999
$ws.Cells.Item(56, 11).Value = 10000
$ws.Cells.Item(56, 12).Value = 49999
$ws.Cells.Item(56, 13).Value = -9286
$ws.Cells.Item(56, 14).Value = -51427

$ws.Cells.Item(92, 8).Value = 0
$ws.Cells.Item(92, 10).Value = 0
$ws.Cells.Item(92, 12).Value = 0
$ws.Cells.Item(92, 14).ClearContents()

$ws.Cells.Item(103, 8).Value = 35899.6
$ws.Cells.Item(103, 10).Value = 35899.6
$ws.Cells.Item(103, 12).Value = 35899.6
$ws.Cells.Item(103, 14).Value = -38243.6

$ws.Cells.Item(139, 8).Value = 200000
$ws.Cells.Item(139, 10).Value = 200000
$ws.Cells.Item(139, 12).Value = 200000
$ws.Cells.Item(139, 14).Value = -210280
